# Update "想去人数" (F column) counts for both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1114
    5  = 47
    8  = 2044
    9  = 7631
    10 = 908
    11 = 426
    13 = 137
    15 = 155
    16 = 7778
    18 = 1347
    23 = 311
    29 = 414
    30 = 623
    31 = 55
    33 = 61
    35 = 38
    36 = 75
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
